# Update the "Data" sheet header labels in J1/K1: the "Student Place of
# Issue" and "Student Date Expiry" headers lose their red "*" (required)
# marker, and scroll the view over to show those columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("K1").Value = "Student Date Expiry"
$ws.Range("J1").Value = "Student Place of Issue"

$ws.Activate()
$ws.Range("G1").Select()
$excel.ActiveWindow.ScrollColumn = 7
